$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.097.35'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '2.518.29'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.08%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.562'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.93%  '
$ws.Range("D9").Value = '2.520.11'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.161'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("E12").Value = '  -2.42%  '
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '2.961.13'
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.74%  '
$ws.Range("D16").Value = '59.016.30'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = '2.516.44'
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.425'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.30%  '
$ws.Range("E26").Value = '  +1.04%  '
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.12%  '
$ws.Range("D30").Value = '0.0₃0772'
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("E31").Value = '  -2.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.15%  '
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("E34").Value = '  -3.51%  '
$ws.Range("E35").Value = '  -8.04%  '
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("E40").Value = '  -1.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '279.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.596'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0931'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("E51").Value = '  -2.00%  '
